$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.814.56"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "1.902.56"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "312.74"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D7").Value = "0.5014"
$ws.Range("E7").Value = "  +3.78%  "
$ws.Range("D8").Value = "0.3809"
$ws.Range("D9").Value = "0.07271"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").Value = "0.9084"
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("D11").Value = "20.82"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "0.07648"
$ws.Range("D13").Value = "1.885.24"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "6.596"
$ws.Range("D16").Value = "91.27"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "0.000008703"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "27.848.32"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").Value = "5.158"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").Value = "154.09"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").Value = "1.866"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").Value = "2.233"
$ws.Range("E26").Value = "  +5.76%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Value = "115.24"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").Value = "4.908"
$ws.Range("D30").Value = "0.08967"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "3.206"
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("D32").Value = "1.231"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7650"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "4.639"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").Value = "0.02056"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "2.544"
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("D37").Value = "1.097"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "0.5554"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "3.014"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").Value = "6.976"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").Value = "8.481"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1510"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.30"
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("D45").Value = "10.59"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "0.4786"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("D49").Value = "67.31"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D51").Value = "0.8989"
$ws.Range("E51").Value = "  -0.48%  "
